{"js": "// Replace the division-problem text runs in the table with the new values.\n// Mapping is a strict 1:1, in-document-order, old-text -> new-text list\n// (old texts are unique across the document, so a targeted search/replace\n// is safe and unambiguous).\nconst replacements = [\n  [\"32\u00f73=\", \"42\u00f72=\"],\n  [\"44\u00f76=\", \"95\u00f72=\"],\n  [\"67\u00f78=\", \"28\u00f72=\"],\n  [\"25\u00f75=\", \"19\u00f76=\"],\n  [\"60\u00f73=\", \"18\u00f72=\"],\n  [\"11\u00f74=\", \"21\u00f72=\"],\n  [\"62\u00f75=\", \"17\u00f73=\"],\n  [\"62\u00f73=\", \"92\u00f74=\"],\n  [\"29\u00f76=\", \"10\u00f73=\"],\n  [\"68\u00f74=\", \"39\u00f79=\"],\n  [\"94\u00f79=\", \"87\u00f79=\"],\n  [\"32\u00f74=\", \"25\u00f73=\"],\n  [\"82\u00f79=\", \"39\u00f79=\"],\n  [\"91\u00f74=\", \"19\u00f73=\"],\n  [\"66\u00f76=\", \"94\u00f74=\"],\n  [\"31\u00f72=\", \"65\u00f74=\"],\n  [\"77\u00f79=\", \"80\u00f79=\"],\n  [\"93\u00f72=\", \"30\u00f73=\"],\n  [\"11\u00f77=\", \"10\u00f79=\"],\n  [\"97\u00f75=\", \"70\u00f72=\"],\n  [\"85\u00f74=\", \"97\u00f75=\"],\n  [\"11\u00f78=\", \"79\u00f78=\"],\n  [\"82\u00f77=\", \"85\u00f73=\"],\n  [\"78\u00f73=\", \"42\u00f75=\"],\n  [\"60\u00f77=\", \"75\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text runs in the table with the new values.\n# Mapping is a strict 1:1 old-text -> new-text list; each old text is\n# unique across the document, so Find/Replace (ReplaceAll) is safe and\n# unambiguous.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"32\u00f73=\", \"42\u00f72=\"),\n    @(\"44\u00f76=\", \"95\u00f72=\"),\n    @(\"67\u00f78=\", \"28\u00f72=\"),\n    @(\"25\u00f75=\", \"19\u00f76=\"),\n    @(\"60\u00f73=\", \"18\u00f72=\"),\n    @(\"11\u00f74=\", \"21\u00f72=\"),\n    @(\"62\u00f75=\", \"17\u00f73=\"),\n    @(\"62\u00f73=\", \"92\u00f74=\"),\n    @(\"29\u00f76=\", \"10\u00f73=\"),\n    @(\"68\u00f74=\", \"39\u00f79=\"),\n    @(\"94\u00f79=\", \"87\u00f79=\"),\n    @(\"32\u00f74=\", \"25\u00f73=\"),\n    @(\"82\u00f79=\", \"39\u00f79=\"),\n    @(\"91\u00f74=\", \"19\u00f73=\"),\n    @(\"66\u00f76=\", \"94\u00f74=\"),\n    @(\"31\u00f72=\", \"65\u00f74=\"),\n    @(\"77\u00f79=\", \"80\u00f79=\"),\n    @(\"93\u00f72=\", \"30\u00f73=\"),\n    @(\"11\u00f77=\", \"10\u00f79=\"),\n    @(\"97\u00f75=\", \"70\u00f72=\"),\n    @(\"85\u00f74=\", \"97\u00f75=\"),\n    @(\"11\u00f78=\", \"79\u00f78=\"),\n    @(\"82\u00f77=\", \"85\u00f73=\"),\n    @(\"78\u00f73=\", \"42\u00f75=\"),\n    @(\"60\u00f77=\", \"75\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
